# unify the conception of DataNode, DataTable, Entity.
# Rename the two worksheets to reflect the new naming scheme, and make the
# second sheet ("DataTable") the active/selected sheet, as in the target
# workbook.

$wb = $excel.ActiveWorkbook

$wsNode  = $wb.Worksheets.Item(1)
$wsTable = $wb.Worksheets.Item(2)

$wsNode.Name  = "DataNode"
$wsTable.Name = "DataTable"

# Switch the active tab from the first sheet to the second one.
$wsTable.Activate()
